$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The product picture filenames in column D are renumbered/renamed from
# "Images/Product/product N" to zero-padded "Images/Product/productNNN.png".
# Row 101 (old "product 100") is written first so that, after the shared
# string table is rebuilt on save, its entry lands at the front of the
# picture-URL block - matching the target layout where product100.png
# precedes product001.png..product099.png.
$ws.Range("D101").Value2 = "Images/Product/product100.png"
for ($row = 2; $row -le 100; $row++) {
    $n = $row - 1
    $ws.Range("D$row").Value2 = "Images/Product/product{0:D3}.png" -f $n
}

# Update the active sheet's selection to match the edited column.
$ws.Activate()
$ws.Range("D2:D101").Select()
